$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 9133.7900000000009
$ws.Range("D10").Value = 7666.84
$ws.Range("E10").Value = 3141.45
$ws.Range("F10").Value = 385

$ws.Range("C11").Value = 3646.37
$ws.Range("D11").Value = 2850.17
$ws.Range("E11").Value = 912.3
$ws.Range("F11").Value = 114

$ws.Range("F14").Select()
